$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "jdxsd"
$ws.Range("B1").Value = "dsbb"
$ws.Range("B1").Select()
